$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so that
# numeric-looking strings (e.g. "1.008", "0.00001099") are not
# auto-converted to numbers when assigned via .Value.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.022.06'
$ws.Range("E2").Value = '  -4.18%  '
$ws.Range("D3").Value = '1.962.48'
$ws.Range("E3").Value = '  -6.13%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").Value = '326.86'
$ws.Range("E5").Value = '  -4.32%  '
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D7").Value = '0.4997'
$ws.Range("E7").Value = '  -5.93%  '
$ws.Range("D8").Value = '0.4202'
$ws.Range("E8").Value = '  -4.24%  '
$ws.Range("D9").Value = '52.88'
$ws.Range("E9").Value = '  -2.99%  '
$ws.Range("D10").Value = '0.09188'
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("D11").Value = '1.097'
$ws.Range("E11").Value = '  -6.63%  '
$ws.Range("D12").Value = '22.89'
$ws.Range("E12").Value = '  -7.29%  '
$ws.Range("D13").Value = '1.964.96'
$ws.Range("E13").Value = '  -5.08%  '
$ws.Range("D14").Value = '7.863'
$ws.Range("E14").Value = '  -8.12%  '
$ws.Range("D15").Value = '6.431'
$ws.Range("E15").Value = '  -6.60%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '0.00001099'
$ws.Range("E17").Value = '  -5.28%  '
$ws.Range("D18").Value = '91.32'
$ws.Range("E18").Value = '  -10.13%  '
$ws.Range("D19").Value = '0.06694'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '19.21'
$ws.Range("E20").Value = '  -9.03%  '
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '5.978'
$ws.Range("E22").Value = '  -5.83%  '
$ws.Range("D23").Value = '29.063.96'
$ws.Range("E23").Value = '  -4.03%  '
$ws.Range("D24").Value = '12.05'
$ws.Range("E24").Value = '  -3.47%  '
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").Value = '2.259.60'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '20.59'
$ws.Range("E27").Value = '  -5.61%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '156.29'
$ws.Range("E28").Value = '  -4.01%  '
$ws.Range("D29").Value = '6.174'
$ws.Range("E29").Value = '  -10.16%  '
$ws.Range("D30").Value = '2.260'
$ws.Range("E30").Value = '  -9.49%  '
$ws.Range("D31").Value = '126.69'
$ws.Range("E31").Value = '  -5.19%  '
$ws.Range("D32").Value = '1.038'
$ws.Range("E32").Value = '  -8.51%  '
$ws.Range("D33").Value = '0.09844'
$ws.Range("E33").Value = '  -6.35%  '
$ws.Range("D34").Value = '1.526'
$ws.Range("E34").Value = '  -8.53%  '
$ws.Range("D35").Value = '5.766'
$ws.Range("E35").Value = '  -8.08%  '
$ws.Range("D36").Value = '3.673'
$ws.Range("E36").Value = '  -6.08%  '
$ws.Range("D37").Value = '0.02416'
$ws.Range("E37").Value = '  -7.82%  '
$ws.Range("D38").Value = '1.296'
$ws.Range("E38").Value = '  -3.32%  '
$ws.Range("D39").Value = '8.940'
$ws.Range("E39").Value = '  -11.69%  '
$ws.Range("D40").Value = '0.06311'
$ws.Range("E40").Value = '  -6.62%  '
$ws.Range("D41").Value = '0.6435'
$ws.Range("E41").Value = '  -7.54%  '
$ws.Range("D42").Value = '11.43'
$ws.Range("E42").Value = '  -9.22%  '
$ws.Range("D43").Value = '0.1983'
$ws.Range("E43").Value = '  -10.29%  '
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = '0.6234'
$ws.Range("E45").Value = '  -7.70%  '
$ws.Range("D46").Value = '13.33'
$ws.Range("E46").Value = '  -6.56%  '
$ws.Range("D47").Value = '2.176'
$ws.Range("E47").Value = '  -8.71%  '
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = '3.460'
$ws.Range("E49").Value = '  -4.75%  '
$ws.Range("D50").Value = '0.00000000333'
$ws.Range("E50").Value = '  -4.77%  '
$ws.Range("D51").Value = '0.06937'
$ws.Range("E51").Value = '  -4.47%  '

# Restore the original (default) cell style now that the text is committed,
# so no stray style index is left attached to these cells.
$ws.Range("D2:E51").Style = "Normal"

